# Add the missing "PM value yesterday" plus a handful of other readings that
# were captured but not yet logged in the raw_data sheet, and extend the
# AM/PM helper-column formula down to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# New readings to append (date/time serials, time-of-day fraction, weight kg)
$newRows = @(
    @{ Row = 162; A = 44097.324999999997; B = 0.32500000000000001; C = 71.5 },
    @{ Row = 163; A = 44097.321527777778; B = 0.3215277777777778;  C = 72.099999999999994 },
    @{ Row = 164; A = 44097.320833333331; B = 0.32083333333333336; C = 72.099999999999994 },
    @{ Row = 165; A = 44097.320138888892; B = 0.32013888888888892; C = 72.099999999999994 },
    @{ Row = 166; A = 44096.895833333336; B = 0.89583333333333337; C = 72.099999999999994 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $item.A
    $cellA.NumberFormat = "m/d/yy h:mm"

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $item.B
    $cellB.NumberFormat = "h:mm"

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = $item.C
}

# Extend the AM/PM helper formula down through the newly added rows
# (D159:D161 already carry the formula; fill it down to D166).
$ws.Range("D162:D166").Formula = "=IF(B162<TIME(12,0,0), ""AM"", ""PM"")"

# Update the view so the newly-entered row is visible/selected, matching
# where the user was working when they made the edit.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 156
$ws.Range("A166").Select() | Out-Null
